$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rearrange the 9 chart images from a single column (rows 0,25,50,...)
# --- into a 3x3 grid (columns 0/10/20, rows 0/20/40), matching the target
# --- drawing1.xml anchors. Positions are expressed in points: Excel COM
# --- Shape.Left/Top are points, and one (0-indexed) grid column/row here
# --- equals the sheet's default column width / row height.

$colWidthPt = $ws.Columns.Item(1).Width
$rowHeightPt = $ws.Rows.Item(1).Height

# shape index (1-based, matches Image 1..Image 9) -> target (col, row) in
# the same 0-indexed col/row units used by the OOXML oneCellAnchor <from>.
$targets = @(
    @{ Shape = 1; Col = 0;  Row = 0  },
    @{ Shape = 2; Col = 10; Row = 0  },
    @{ Shape = 3; Col = 20; Row = 0  },
    @{ Shape = 4; Col = 0;  Row = 20 },
    @{ Shape = 5; Col = 10; Row = 20 },
    @{ Shape = 6; Col = 20; Row = 20 },
    @{ Shape = 7; Col = 0;  Row = 40 },
    @{ Shape = 8; Col = 10; Row = 40 },
    @{ Shape = 9; Col = 20; Row = 40 }
)

foreach ($t in $targets) {
    $shp = $ws.Shapes.Item($t.Shape)
    $shp.Left = $t.Col * $colWidthPt
    $shp.Top = $t.Row * $rowHeightPt
}

# --- The sheet's used range grows to A1:U41 (the new grid's bounding box,
# --- 0-indexed col 20 / row 40 -> 1-indexed column U / row 41), and
# --- sheetData picks up placeholder rows at r=1, r=21 and r=41.
#
# Touch the two extreme corner cells so the sheet dimension expands to
# A1:U41, then clear any formatting trace they leave behind so they stay
# empty cells.
$corner1 = $ws.Cells.Item(1, 1)
$corner1.NumberFormat = "General"
$corner1.ClearFormats()

$corner2 = $ws.Cells.Item(41, 21)
$corner2.NumberFormat = "General"
$corner2.ClearFormats()

# Row 21 picks up an empty <row r="21"/> stub without any cell content;
# toggling OutlineLevel touches the row without creating a cell.
$ws.Rows.Item(21).OutlineLevel = 1
$ws.Rows.Item(21).OutlineLevel = 0
